$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "TOR130", "2024", 45545.875, "Tuesday")
    ,@(3, "TOR130", "2023", 45181.875, "Tuesday")
    ,@(4, "TOR130", "2022", 44817.875, "Tuesday")
    ,@(5, "TOR130", "2019", 43718.875, "Tuesday")
    ,@(6, "TOR130", "2018", 43354.875, "Tuesday")
    ,@(7, "TOR130", "2017", 42991.875, "Wednesday")
    ,@(8, "TOR130", "2021", 44453.875, "Tuesday")
    ,@(9, "TOR330", "2023", 45179.41666666666, "Sunday")
    ,@(10, "TOR330", "2021", 44451.41666666666, "Sunday")
    ,@(11, "TOR330", "2017", 42988.41666666666, "Sunday")
    ,@(12, "TOR330", "2024", 45543.41666666666, "Sunday")
    ,@(13, "TOR330", "2013", 41525.41666666666, "Sunday")
    ,@(14, "TOR330", "2022", 44815.41666666666, "Sunday")
    ,@(15, "TOR330", "2014", 41889.41666666666, "Sunday")
    ,@(16, "TOR330", "2019", 43716.41666666666, "Sunday")
    ,@(17, "TOR330", "2018", 43352.41666666666, "Sunday")
    ,@(18, "TOR330", "2016", 42624.41666666666, "Sunday")
    ,@(19, "TOR330", "2012", 41161.41666666666, "Sunday")
    ,@(20, "TOR330", "2011", 40797.41666666666, "Sunday")
    ,@(21, "TOR330", "2015", 42260.41666666666, "Sunday")
    ,@(22, "TOR330", "2010", 40433.41666666666, "Sunday")
    ,@(23, "TOR450", "2023", 45177.83333333334, "Friday")
    ,@(24, "TOR450", "2024", 45541.83333333334, "Friday")
    ,@(25, "TOR450", "2022", 44813.83333333334, "Friday")
    ,@(26, "TOR450", "2019", 43714.83333333334, "Friday")
    ,@(27, "TOR450", "2021", 44449.83333333334, "Friday")
)

foreach ($row in $data) {
    $r = $row[0]
    $race = $row[1]
    $year = $row[2]
    $cdate = $row[3]
    $dow = $row[4]

    $ws.Range("A$r").Value = $race

    $cellB = $ws.Range("B$r")
    $cellB.NumberFormat = "@"
    $cellB.Value = $year
    $cellB.Style = "Normal"

    $cellC = $ws.Range("C$r")
    $cellC.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cellC.Value = $cdate

    $ws.Range("D$r").Value = $dow
}

Write-Host "Done. Dimension: $($ws.UsedRange.Address())"
